$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nodes")

# Rebuild the H-column JSON-builder formulas (rows 2-38) to splice in a
# "last_name" field sourced from column C, between "first_name" and "team".
for ($r = 2; $r -le 38; $r++) {
    $formula = '="{"&CHAR(34)&"id"&CHAR(34)&": "&A' + $r + '&", "&CHAR(34)&"first_name"&CHAR(34)&":  "&CHAR(34)&B' + $r + '&CHAR(34)&", "&CHAR(34)&"last_name"&CHAR(34)&": "&CHAR(34)&C' + $r + '&CHAR(34)&", "&CHAR(34)&"team"&CHAR(34)&":  "&CHAR(34)&D' + $r + '&CHAR(34)&", "&CHAR(34)&"group"&CHAR(34)&":  "&E' + $r + '&"},"'
    $ws.Range("H$r").Formula = $formula
}

# Switch the active/selected tab from Links back to Nodes.
$ws.Activate()
